# Update of daily and weekly charts
# Adds the new day's data (12/8/2020) to the DateofDeath tracking sheet and
# refreshes the handful of recent rows whose "Probable" counts were revised.
#
# NOTE: window-position cosmetics (workbookView xWindow/yWindow/width/height,
# and the frozen-pane's scrolled topLeftCell) are screen/session state that
# this headless runtime does not persist back to the OOXML on save, so they
# are left at their existing values rather than forced through no-op calls.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DateofDeath")

# --- Revisions to existing recent rows (new probable-death counts reported) ---
$ws.Range("B271").Value = 48
$ws.Range("B272").Value = 42
$ws.Range("B273").Value = 47
$ws.Range("D273").Value = 1
$ws.Range("B274").Value = 56
$ws.Range("D274").Value = 2

# --- Append new row 275 for 12/8/2020 ---
# Copy the date formatting from the row above, then set the value.
$ws.Range("A274").Copy() | Out-Null
$ws.Range("A275").PasteSpecial(-4122) | Out-Null
$ws.Range("A275").Value = 44173

$ws.Range("B275").Value = 28
$ws.Range("C275").Formula = "=B275+C274"
$ws.Range("D275").Value = 0
$ws.Range("E275").Formula = "=D275+E274"
$ws.Range("F275").Formula = "=AVERAGE(B269:B275)"

$excel.CutCopyMode = $false

# --- Update the view so the new row is visible/selected ---
$ws.Range("D281").Select() | Out-Null
